# feat: add 2022-Q1 data
#
# The existing "总计" (Total) sheet is repurposed in-place to become the new
# "2022-Q1" quarterly holdings sheet (same underlying sheet/tab position,
# just renamed + its content replaced with the new fund-holding table), and
# a brand-new "总计" sheet is appended right after it, carrying the refreshed
# totals table (the previous totals, shifted down one row, with a new
# 2022-Q1 row inserted at the top).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Structural change: rename "总计" -> "2022-Q1", then add a fresh
#    "总计" sheet right after it (so tab order stays
#    ... 2021-Q4, 2022-Q1, 总计). The new sheet is created by duplicating
#    an existing, well-formed data sheet (same sheetPr/margins template
#    as the original "总计" sheet) so its look & feel matches the rest of
#    the workbook instead of a bare blank sheet.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

$template = $wb.Worksheets.Item("2021-Q4")
$template.Copy($null, $q1)
$total = $wb.Worksheets.Item("2021-Q4 (2)")
$total.Name = "总计"

# ---------------------------------------------------------------------
# 2. Rebuild "2022-Q1" sheet content (was the "总计" sheet's sheet file).
#    Old data lived in A1:D6 - clear rows below the header/first data row,
#    then widen the header style out to the new columns E:H, and write
#    the new fund-holding table.
# ---------------------------------------------------------------------
$q1.Range("A3:D6").Clear()

# Extend header styling (s="2" incl. border/font) from D1 across E1:H1.
$q1.Range("D1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1.Range("A2").Value = 0

# Numeric-looking text fields must stay text (keep the trailing zeros /
# leading digits exactly as authored) instead of being auto-coerced to
# numbers. Format as Text first so the literal string sticks, then drop
# back to the Normal style (keeps the stored value a string without
# leaving a quote-prefixed/text-number-format style behind on the cell).
$numericTextCells = "B2", "D2", "E2", "F2", "G2"
foreach ($addr in $numericTextCells) {
    $q1.Range($addr).NumberFormat = "@"
}
$q1.Range("B2").Value = "202801"
$q1.Range("C2").Value = "南方全球精选配置(QDII-FOF)"
$q1.Range("D2").Value = "18.00"
$q1.Range("E2").Value = "28.82"
$q1.Range("F2").Value = "1.80"
$q1.Range("G2").Value = "0.3240"
foreach ($addr in $numericTextCells) {
    $q1.Range($addr).Style = "Normal"
}
$q1.Range("H2").Value = 2

# ---------------------------------------------------------------------
# 3. Populate the brand-new "总计" sheet with the refreshed totals table.
#    It currently still holds the copied "2021-Q4" fund-holding data
#    (A1:H2) - clear everything first, then re-apply the totals-table
#    header/row styling (copied from the old "总计"/now "2022-Q1" sheet,
#    which still uses the correct style) and write the values.
# ---------------------------------------------------------------------
$total.Cells.Clear()

$q1.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)

$q1.Range("A2").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.32

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.38

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 1
$total.Range("D4").Value = 0.43

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 2
$total.Range("D5").Value = 0.79

$total.Range("A6").Value = 4
$total.Range("B6").Value = "2021-Q1"
$total.Range("C6").Value = 4
$total.Range("D6").Value = 0.99

$total.Range("A7").Value = 5
$total.Range("B7").Value = "2020-Q4"
$total.Range("C7").Value = 2
$total.Range("D7").Value = 1.01

# ---------------------------------------------------------------------
# 4. Restore the original active sheet/selection (2020-Q4, first tab).
# ---------------------------------------------------------------------
$wb.Worksheets.Item(1).Activate()
